# Inserts one new weekly price record for "Coliflor" (Terminal La Palmera de
# La Serena) above the existing row 967, pushing every subsequent row down
# by one. This mirrors the author's "Fruta / hortaliza, semanal" update,
# which adds a new week's worth of data to the bottom of the logical
# sequence (the sheet stores rows oldest-last, so the new pair of
# Primera/Segunda quality rows lands just above the row that used to be
# first in that pair sequence).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 967; everything from the old row 967
# downward (through 1044) shifts down to 968-1045.
$ws.Rows.Item(967).Insert()

# Populate the newly inserted row 967 with the new weekly record.
$ws.Range("A967").Value = 8
$ws.Range("B967").Value = "Terminal La Palmera de La Serena"
$ws.Range("C967").Value = "Coquimbo"
$ws.Range("D967").Value = 45013
$ws.Range("E967").Value = 4
$ws.Range("F967").Value = 100112008
$ws.Range("G967").Value = "Coliflor"
$ws.Range("H967").Value = "Sin especificar"
$ws.Range("I967").Value = "Segunda"
$ws.Range("J967").Value = 1500
$ws.Range("K967").Value = 900
$ws.Range("L967").Value = 1000
$ws.Range("M967").Value = 950
$ws.Range("N967").Value = "$/unidad"
$ws.Range("O967").Value = "Provincia del Elquí"
$ws.Range("P967").Value = 950
$ws.Range("Q967").Value = 1
$ws.Range("R967").Value = "Hortaliza"
